$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each crypto row with refreshed data.
# For D-column values that look numeric (e.g. "0.9996"), force text format first so
# Excel keeps them as literal strings (matching the original "Price" text formatting)
# instead of auto-converting them into floating point numbers; the style is then reset
# back to Normal so no extra formatting is introduced.

$ws.Range("D2").Value = "29.147.37"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").Value = "1.855.03"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6904"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07702"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3027"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08083"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").Value = "1.838.18"
$ws.Range("E12").Value = "  -0.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7153"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.142"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.35%  "

$ws.Range("D16").Value = "29.144.41"
$ws.Range("E16").Value = "  -0.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.704"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007723"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "235.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9996"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "2.108.65"
$ws.Range("E22").Value = "  +0.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.426"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.980"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("E27").Value = "  -2.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.940"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.405"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.458"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.485"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.985"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05165"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.164"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7037"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9993"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.652"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01843"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.05%  "

$ws.Range("E40").Value = "  +1.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9348"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.72%  "

$ws.Range("D42").Value = "1.109.62"
$ws.Range("E42").Value = "  +5.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4257"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.865"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.29%  "

$ws.Range("E48").Value = "  +2.42%  "

$ws.Range("D49").Value = "2.005.47"
$ws.Range("E49").Value = "  +0.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.114"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.939"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.29%  "
